$d = $word.ActiveDocument
$range = $d.Content
$found = $range.Find.Execute("preferencial). utilizador")
Write-Host "Found: $found start=$($range.Start) end=$($range.End)"

$prefixLen = "preferencial). ".Length
$insStart = $range.Start + $prefixLen
$insRange = $d.Range($insStart, $insStart)
$insRange.Text = "O "
